$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" note text ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.0 = 27415.97 pesos`n✅ 27415.97 pesos = 6.95 = 942.86 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas sheet: update N10/O10 and N12/O12 rate values ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 142.8
$tasas.Range("O10").Value = 3915
$tasas.Range("N12").Value = 3945.95
$tasas.Range("O12").Value = 135.705
